# Insert a new weekly record as row 19, pushing existing rows 19:146 down to 20:147.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(19).Insert()

$ws.Cells.Item(19, 1).Value = 1
$ws.Cells.Item(19, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(19, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(19, 4).Value = 44490
$ws.Cells.Item(19, 4).NumberFormat = $ws.Cells.Item(20, 4).NumberFormat
$ws.Cells.Item(19, 5).Value = 15
$ws.Cells.Item(19, 6).Value = "Fruta"
$ws.Cells.Item(19, 7).Value = 100108
$ws.Cells.Item(19, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(19, 9).Value = 100108006
$ws.Cells.Item(19, 10).Value = "Plátano"
$ws.Cells.Item(19, 11).Value = "Sin especificar"
$ws.Cells.Item(19, 12).Value = "Pintón"
$ws.Cells.Item(19, 13).Value = 120
$ws.Cells.Item(19, 14).Value = 24000
$ws.Cells.Item(19, 15).Value = 25000
$ws.Cells.Item(19, 16).Value = 24500
$ws.Cells.Item(19, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(19, 18).Value = "Ecuador"
$ws.Cells.Item(19, 19).Value = 1225
$ws.Cells.Item(19, 20).Value = 20
